$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "48"
$ws.Range("C9").Characters(27, 10).Text = "11/27/2023"
$ws.Range("C9").Characters(48, 10).Text = "12/3/2023"

# --- Crime statistics table updates (rows 14-29) ---
# Row 14
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 11
$ws.Range("K14").Value = 57.142857142857
$ws.Range("L14").Value = -31.25
$ws.Range("M14").Value = 57.142857142857
$ws.Range("N14").Value = -79.245283018867

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 27
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 8
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -12.903225806451
$ws.Range("N15").Value = -66.25

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 196
$ws.Range("J16").Value = 205
$ws.Range("K16").Value = -4.390243902439
$ws.Range("L16").Value = 18.787878787878
$ws.Range("M16").Value = -26.591760299625
$ws.Range("N16").Value = -76.88679245283

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -54.545454545454
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -48.717948717948
$ws.Range("I17").Value = 422
$ws.Range("J17").Value = 486
$ws.Range("K17").Value = -13.168724279835
$ws.Range("L17").Value = -14.052953156822
$ws.Range("M17").Value = 42.087542087542
$ws.Range("N17").Value = -54.525862068965

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -7.692307692307
$ws.Range("I18").Value = 161
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -13.440860215053
$ws.Range("L18").Value = 1.898734177215
$ws.Range("M18").Value = 37.606837606837
$ws.Range("N18").Value = -82.229580573951

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 266.666666666667
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 36
$ws.Range("I19").Value = 411
$ws.Range("J19").Value = 318
$ws.Range("K19").Value = 29.245283018867
$ws.Range("L19").Value = 35.197368421052
$ws.Range("M19").Value = 45.229681978798
$ws.Range("N19").Value = 7.874015748031

# Row 20
$ws.Range("A14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "'0"
$ws.Range("D16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 200
$ws.Range("J20").Value = 68
$ws.Range("K20").Value = 92.647058823529
$ws.Range("L20").Value = 142.592592592593
$ws.Range("M20").Value = 167.34693877551
$ws.Range("N20").Value = -45.188284518828

# Row 21
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -3.883495145631
$ws.Range("I21").Value = 1359
$ws.Range("J21").Value = 1295
$ws.Range("K21").Value = 4.942084942084
$ws.Range("L21").Value = 12.686567164179
$ws.Range("M21").Value = 29.305423406279
$ws.Range("N21").Value = -60.436681222707

# Row 22
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 15
$ws.Range("K22").Value = -21.052631578947
$ws.Range("L22").Value = -6.25
$ws.Range("M22").Value = 66.666666666666

# Row 23
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -36
$ws.Range("I23").Value = 224
$ws.Range("J23").Value = 220
$ws.Range("K23").Value = 1.818181818181
$ws.Range("L23").Value = 6.161137440758
$ws.Range("M23").Value = 39.130434782608

# Row 24
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 74
$ws.Range("H24").Value = 42.307692307692
$ws.Range("I24").Value = 818
$ws.Range("J24").Value = 764
$ws.Range("K24").Value = 7.068062827225
$ws.Range("L24").Value = 14.726507713885
$ws.Range("M24").Value = 39.115646258503

# Row 25
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 62
$ws.Range("H25").Value = 5.084745762711
$ws.Range("I25").Value = 580
$ws.Range("J25").Value = 566
$ws.Range("K25").Value = 2.473498233215
$ws.Range("L25").Value = 0.519930675909
$ws.Range("M25").Value = -26.582278481012

# Row 26
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 43
$ws.Range("J26").Value = 44
$ws.Range("K26").Value = -2.272727272727
$ws.Range("L26").Value = 38.709677419354

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 51
$ws.Range("J27").Value = 67
$ws.Range("K27").Value = -23.880597014925
$ws.Range("L27").Value = -21.538461538461

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = -17.142857142857
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -23.684210526315
$ws.Range("N28").Value = -76.8

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 26
$ws.Range("J29").Value = 29
$ws.Range("K29").Value = -10.344827586206
$ws.Range("L29").Value = -44.680851063829
$ws.Range("M29").Value = -18.75
$ws.Range("N29").Value = -76.991150442477
